# Updated cryptos list (GitHub Actions scheduled refresh).
# Price/volume text cells are stored as plain text (t="inlineStr"/shared
# string), not numbers, so for any new value that looks numeric we force
# the cell to Text format first (otherwise Excel would silently convert
# e.g. "1.00" -> 1 or "0.170" -> 0.17) and then restore the "Normal"
# style afterwards so no stray number-format style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.783.59'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = '3.379.30'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.379.00'
$ws.Range('E8').Value = '  -1.10%  '
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.56'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('D13').Value = '3.959.52'
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.25%  '
$ws.Range('E16').Value = '  -3.63%  '
$ws.Range('D17').Value = '3.379.28'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '60.901.80'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '376.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.62%  '
$ws.Range('D24').Value = '3.526.74'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000125'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.170'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.72'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('E36').Value = '  -4.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.85'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.48'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('E40').Value = '  -4.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.771'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.428.28'
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.22%  '
